$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> index map for the columns that change
$colIndex = @{ 'D' = 4; 'J' = 10; 'K' = 11; 'L' = 12; 'M' = 13; 'P' = 16 }

$changes = @{
    2  = @{ D = 44895; J = 200;  K = 1200; L = 1300; M = 1255; P = 1255 }
    3  = @{ D = 44638; J = 800;  K = 2500; L = 2800; M = 2650; P = 2650 }
    4  = @{ D = 44537; J = 800;  K = 1300; L = 1400; M = 1350; P = 1350 }
    5  = @{ D = 44210; J = 1450; K = 1600; L = 1700; M = 1650; P = 1650 }
    6  = @{ D = 44893; J = 3300; K = 1200; L = 1300; M = 1261; P = 1261 }
    7  = @{ D = 44907; J = 2300; K = 900;  L = 1000; M = 952;  P = 952  }
    8  = @{ D = 44883; J = 290;  K = 1400; L = 1500; M = 1434; P = 1434 }
    9  = @{ D = 44200; J = 1500; K = 1400; L = 1500; M = 1450; P = 1450 }
    10 = @{ D = 44175; J = 1400; K = 1900; L = 2000; M = 1950; P = 1950 }
}

foreach ($row in $changes.Keys) {
    $rowChanges = $changes[$row]
    foreach ($col in $rowChanges.Keys) {
        $c = $colIndex[$col]
        $ws.Cells.Item($row, $c).Value2 = $rowChanges[$col]
    }
}
